# Updates cryptos list values (Price / Volume(1h) columns) per the
# upstream GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.405.04'
$ws.Range("E2").Value = '  +6.73%  '
$ws.Range("D3").Value = '3.572.27'
$ws.Range("E3").Value = '  +2.38%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = "'416.97"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").Value = "'129.94"
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("D7").Value = "'0.648"
$ws.Range("E7").Value = '  +2.79%  '
$ws.Range("D8").Value = '3.562.72'
$ws.Range("E8").Value = '  +2.30%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = "'0.764"
$ws.Range("E10").Value = '  +3.90%  '
$ws.Range("D11").Value = "'0.176"
$ws.Range("E11").Value = '  +13.26%  '
$ws.Range("E12").Value = '  +49.12%  '
$ws.Range("D13").Value = "'42.30"
$ws.Range("E13").Value = '  -0.92%  '
$ws.Range("D14").Value = "'9.99"
$ws.Range("E14").Value = '  +1.74%  '
$ws.Range("D15").Value = '4.140.59'
$ws.Range("E15").Value = '  +2.51%  '
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").Value = "'20.42"
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("D18").Value = '3.576.74'
$ws.Range("E18").Value = '  +3.59%  '
$ws.Range("E19").Value = '  +5.00%  '
$ws.Range("D20").Value = '67.341.38'
$ws.Range("E20").Value = '  +6.84%  '
$ws.Range("D21").Value = "'12.29"
$ws.Range("E21").Value = '  -3.42%  '
$ws.Range("D22").Value = "'461.65"
$ws.Range("E22").Value = '  -1.32%  '
$ws.Range("D23").Value = "'88.01"
$ws.Range("E23").Value = '  -3.18%  '
$ws.Range("E24").Value = '  -5.97%  '
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("D26").Value = "'3.35"
$ws.Range("E26").Value = '  +0.85%  '
$ws.Range("D27").Value = "'10.19"
$ws.Range("E27").Value = '  -5.14%  '
$ws.Range("D28").Value = "'35.31"
$ws.Range("E28").Value = '  +4.50%  '
$ws.Range("E29").Value = '  +0.41%  '
$ws.Range("D30").Value = "'2.79"
$ws.Range("E30").Value = '  +4.12%  '
$ws.Range("E31").Value = '  +1.96%  '
$ws.Range("D32").Value = "'7.42"
$ws.Range("E32").Value = '  -2.02%  '
$ws.Range("E33").Value = '  +3.76%  '
$ws.Range("D34").Value = "'41.61"
$ws.Range("E34").Value = '  +1.15%  '
$ws.Range("E35").Value = '  -4.46%  '
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").Value = "'56.67"
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("D38").Value = "'0.0493"
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("D39").Value = '0.0₃0704'
$ws.Range("E39").Value = '  +18.92%  '
$ws.Range("D40").Value = "'0.146"
$ws.Range("E40").Value = '  +8.17%  '
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("D42").Value = "'3.03"
$ws.Range("E42").Value = '  -1.06%  '
$ws.Range("D43").Value = "'148.36"
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("E45").Value = '  -2.36%  '
$ws.Range("E46").Value = '  -2.90%  '
$ws.Range("D47").Value = "'0.310"
$ws.Range("E47").Value = '  -4.06%  '
$ws.Range("D48").Value = "'1.97"
$ws.Range("E48").Value = '  -4.72%  '
$ws.Range("E49").Value = '  -1.58%  '
$ws.Range("D50").Value = "'2.70"
$ws.Range("E50").Value = '  +15.95%  '
$ws.Range("D51").Value = "'15.62"
$ws.Range("E51").Value = '  -4.97%  '

Write-Output "Updated 84 cells"
